# Update the "ejemplo_carga_masiva_postman.xlsx" example workbook used
# for Postman bulk-load testing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings are appended to sst in the order they are first
# encountered while saving. To reproduce the exact ordering seen in the
# target workbook (all three new "id" values, then all three new "codigo"
# values), write column A for all rows first, then column C for all rows,
# before touching the remaining (already-interned) columns.

# --- Column A (id) ---
$ws.Range("A2").Value = "PROD_POSTMAN_004"
$ws.Range("A3").Value = "PROD_POSTMAN_005"
$ws.Range("A4").Value = "PROD_POSTMAN_006"

# --- Column C (codigo) ---
$ws.Range("C2").Value = "AMX500-POSTMAN2"
$ws.Range("C3").Value = "IBU400-POSTMAN2"
$ws.Range("C4").Value = "LOS50-POSTMAN2"

# --- Row 2: Amoxicilina 500mg ---
$ws.Range("B2").Value = "Amoxicilina 500mg"
$ws.Range("D2").Value = "ANTIBIOTICS"
$ws.Range("E2").Value = "Cápsula"
$ws.Range("F2").Value = 1550
$ws.Range("G2").Value = "CERT-INVIMA-2024-001"
$ws.Range("H2").Value = "Temperatura ambiente 15-30°C"
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 30

# --- Row 3: Ibuprofeno 400mg ---
$ws.Range("B3").Value = "Ibuprofeno 400mg"
$ws.Range("D3").Value = "ANALGESICS"
$ws.Range("E3").Value = "Tableta"
$ws.Range("F3").Value = 850
$ws.Range("G3").Value = "CERT-INVIMA-2024-002"
$ws.Range("H3").Value = "Temperatura ambiente"
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 170
$ws.Range("K3").Value = 50

# --- Row 4: Losartán 50mg ---
$ws.Range("B4").Value = "Losartán 50mg"
$ws.Range("D4").Value = "CARDIOVASCULARES"
$ws.Range("E4").Value = "Tableta"
$ws.Range("F4").Value = 1200
$ws.Range("G4").Value = "CERT-INVIMA-2024-003"
$ws.Range("H4").Value = "Temperatura ambiente 15-25°C"
$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 40

# --- Column widths (nearest values achievable through the ColumnWidth
#     property's internal quantization that reproduce the stored widths) ---
$ws.Columns.Item(3).ColumnWidth = 24.333333333333332
$ws.Columns.Item(4).ColumnWidth = 14.166666666666666
$ws.Columns.Item(9).ColumnWidth = 19.5
$ws.Columns.Item(10).ColumnWidth = 18.5

# --- Selection ---
$ws.Range("K7").Select()

$wb.Save()
